$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Reserva-0001")
$after = $wb.Worksheets.Item("Reserva-0002")

$source.Copy($null, $after)

$newSheet = $wb.Worksheets.Item($after.Index + 1)
$newSheet.Name = "Residente-0001"
